# Update the "Các thành phần cần thay đổi" sheet: shorten the assignee
# names "Giang" -> "G" and "Dinh" -> "D", and clear the two price/budget
# notes ("1tr" / "3tr") in the small task table further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Các thành phần cần thay đổi")

$ws.Range("E5").Value = "D"
$ws.Range("E6").Value = "D"

$ws.Range("E2").Value = "G"
$ws.Range("E4").Value = "G"
$ws.Range("E7").Value = "G"

$ws.Range("D18").ClearContents()
$ws.Range("D20").ClearContents()

$ws.Range("G7").Select()
